$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108
$r1.VerticalAlignment = -4160
$r1.Borders.LineStyle = 1
$r1.Borders.Weight = 2

$r1.Copy()
$r2 = $ws.Range("A2")
$r2.PasteSpecial(-4122)
$excel.CutCopyMode = $false
